$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking strings (e.g. "578.04") as
# plain text. Mark each touched price cell as Text first so Excel does
# not silently coerce the assigned value into a real number (which would
# drop formatting such as trailing zeros).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Row 27 and Row 28 swap places (RenderToken overtakes Dai) and get new
# price / volume figures.
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  -1.99%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.02%  "

# Refreshed price / volume(1h) figures for the remaining rows.
$ws.Range("D2").Value = '66.862.05'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '3.112.30'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '578.04'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").Value = '172.53'
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.107.22'
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value = '  -0.86%  '
$ws.Range("D10").Value = '6.45'
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("D11").Value = '0.153'
$ws.Range("E11").Value = '  -1.58%  '
$ws.Range("D12").Value = '0.481'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").Value = '0.0000246'
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("D14").Value = '37.43'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("D16").Value = '3.623.94'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '66.770.42'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '7.14'
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").Value = '3.107.99'
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '16.44'
$ws.Range("E20").Value = '  +1.98%  '
$ws.Range("D21").Value = '478.05'
$ws.Range("E21").Value = '  +1.42%  '
$ws.Range("D22").Value = '0.715'
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = '7.99'
$ws.Range("E23").Value = '  +5.90%  '
$ws.Range("D24").Value = '13.55'
$ws.Range("E24").Value = '  +5.18%  '
$ws.Range("D25").Value = '83.97'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D29").Value = '2.44'
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").Value = '7.94'
$ws.Range("E30").Value = '  -3.07%  '
$ws.Range("E31").Value = '  -1.09%  '
$ws.Range("D32").Value = '28.66'
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("D34").Value = '0.0₃0943'
$ws.Range("E34").Value = '  -7.48%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = '5.86'
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").Value = '0.978'
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("D38").Value = '47.29'
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("D39").Value = '2.09'
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("D40").Value = '50.04'
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("D41").Value = '0.310'
$ws.Range("E41").Value = '  -2.26%  '
$ws.Range("E42").Value = '  -1.78%  '
$ws.Range("D43").Value = '8.66'
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("D44").Value = '2.806.91'
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").Value = '0.0355'
$ws.Range("E45").Value = '  -2.67%  '
$ws.Range("D46").Value = '378.91'
$ws.Range("E46").Value = '  -4.58%  '
$ws.Range("D47").Value = '2.56'
$ws.Range("E47").Value = '  -11.10%  '
$ws.Range("D48").Value = '136.18'
$ws.Range("E48").Value = '  +0.67%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").Value = '25.02'
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("E51").Value = '  -2.12%  '
